# Ozeret-mode sign-in button fix
#
# Two visitors (Visitor 1 / Visitor 2) sign in on the "06-07-2022" daily
# sheet: their "Time In" cell moves from the red "Not in Bunk" placeholder
# to an actual check-in time, highlighted green (the same red/green
# fill pair already used elsewhere in the workbook for sign-in state).
# The curfew times for that night also get updated, and the Key sheet's
# tally columns (On Time / Late) pick up the resulting counts, which in
# turn feed the "# Returned" / "# Still Out" counters on the daily sheet.

$wb = $excel.ActiveWorkbook

# --- "06-07-2022" daily attendance sheet -----------------------------
$today = $wb.Worksheets.Item("06-07-2022")

# Curfew times for the night (Leaving Camp Curfew / Night Off Curfew rows)
$today.Range("I2").Value2 = "9:55 PM"
$today.Range("I4").Value2 = "9:57 PM"

# Visitor 2 (row 2) and Visitor 1 (row 4) sign back in - the ozeret
# "sign in" button writes the actual time and flips the cell to the
# green "signed in" fill (previously stuck on the red "Not in Bunk"
# placeholder because the button didn't work in ozeret mode).
$today.Range("E2").Value2 = "9:54 PM"
$today.Range("E2").Interior.Color = 13492663   # RGB(183,225,205) / B7E1CD, BGR-packed
$today.Range("E4").Value2 = "9:54 PM"
$today.Range("E4").Interior.Color = 13492663   # RGB(183,225,205) / B7E1CD, BGR-packed

# With both visitors back, 2 more have returned and 2 fewer are still out
$today.Range("I7").Value2 = 2
$today.Range("I8").Value2 = 4

# --- "Key" sheet: per-person On Time / Late tallies -------------------
$key = $wb.Worksheets.Item("Key")

$key.Range("E2").Value2 = 1   # Staff Member 1 - Late
$key.Range("D3").Value2 = 1   # Staff Member 2 - On Time
$key.Range("E4").Value2 = 1   # Staff Member 3 - Late
$key.Range("E5").Value2 = 1   # Staff Member 4 - Late
$key.Range("D6").Value2 = 2   # Visitor 1 - On Time
$key.Range("D7").Value2 = 2   # Visitor 2 - On Time
